# Apply the lookups.xlsx changes:
#  - Companies:  add header "active" in B1 (copy of A1 style), set B2 from FALSE -> TRUE
#  - Locations:  add header "company" in B1 (copy of A1 style) and several new rows
#  - AssetTypes: add header "location" in B1 (copy of A1 style) and several new rows

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: Companies
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Companies")

# New header cell B1 ("active"), matching the bold/bordered header style used by A1.
$ws1.Range("A1").Copy()
$ws1.Range("B1").PasteSpecial(-4122)
$ws1.Range("B1").Value = "active"

# B2 flips from FALSE to TRUE (stored as text, like the rest of the sheet,
# not a real boolean). Plain `.Value = "TRUE"` auto-coerces to a boolean
# cell, so build the text via a formula that evaluates to a string and then
# flatten the cell back down to a literal value with copy / paste-values.
$b2 = $ws1.Range("B2")
$b2.Formula = '="TRUE"'
$b2.Copy()
$b2.PasteSpecial(-4163)

# ---------------------------------------------------------------------------
# Sheet 2: Locations
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Locations")

# New header cell B1 ("company"), matching A1's header style.
$ws2.Range("A1").Copy()
$ws2.Range("B1").PasteSpecial(-4122)
$ws2.Range("B1").Value = "company"

# Existing row 2 (A2 = test123) becomes row 4; new rows are inserted/appended
# so the final layout is:
#   2: BC
#   3: AB      | NHS
#   4: test123 | NHS
#   5: mhfkymh
#   6: feqF    | NHS
$ws2.Range("A2").Value = "BC"

$ws2.Range("A3").Value = "AB"
$ws2.Range("B3").Value = "NHS"

$ws2.Range("A4").Value = "test123"
$ws2.Range("B4").Value = "NHS"

$ws2.Range("A5").Value = "mhfkymh"

$ws2.Range("A6").Value = "feqF"
$ws2.Range("B6").Value = "NHS"

# ---------------------------------------------------------------------------
# Sheet 3: AssetTypes
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("AssetTypes")

# New header cell B1 ("location"), matching A1's header style.
$ws3.Range("A1").Copy()
$ws3.Range("B1").PasteSpecial(-4122)
$ws3.Range("B1").Value = "location"

# Row 2 ("banana") is replaced by "cableway" + location "AB"; rows 3-9 are new.
$ws3.Range("A2").Value = "cableway"
$ws3.Range("B2").Value = "AB"

$ws3.Range("A3").Value = "test"

$ws3.Range("A4").Value = "hiuewkcs"

# Stored as text ("123"), not a number, so force it the same way as B2 above.
$a5 = $ws3.Range("A5")
$a5.Formula = '="123"'
$a5.Copy()
$a5.PasteSpecial(-4163)

$ws3.Range("A6").Value = "pppoknl"

$ws3.Range("A7").Value = "okjhgnb"

$ws3.Range("A8").Value = "u54jktyu"
$ws3.Range("B8").Value = "AB"

$ws3.Range("A9").Value = "VDS"
$ws3.Range("B9").Value = "feqF"
